$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("boosters")

# Insert two new columns (pg, vg) before the existing "nicotine" column (C),
# shifting nicotine/volume right into E/F.
$ws.Range("C1:C8").Insert(-4161)
$ws.Range("C1:C8").Insert(-4161)

# New header cells, matching style of the other header cells.
$ws.Range("C1").Value = "pg"
$ws.Range("D1").Value = "vg"
$ws.Range("C1:D1").Style = $ws.Range("B1").Style

# New column widths for the freshly inserted columns (nearest reachable
# quantized ColumnWidth values for target raw widths 13.42578125 / 12.7109375).
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 11.833333333333334

# pg / vg data values.
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 70

$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 70

$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 70

$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 50

$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 50

$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 50

$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 50

$ws.Range("D5").Select()

$wb.Windows.Item(1).Left = 2010
$wb.Windows.Item(1).Top = 7965
